# Weekly update: insert two new price rows at the top of the Pimiento
# data block (row 582) and push the existing history down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 582:624 down to 584:626 by inserting two blank rows.
$ws.Range("A582:R583").Insert()

# New row 582: Zafiro rojo, Primera, 2024-? (serial 45265)
$ws.Range("A582").Value = 7
$ws.Range("B582").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C582").Value = "Ñuble"
$ws.Range("D582").Value = 45265
$ws.Range("E582").Value = 16
$ws.Range("F582").Value = 100112002
$ws.Range("G582").Value = "Pimiento"
$ws.Range("H582").Value = "Zafiro rojo"
$ws.Range("I582").Value = "Primera"
$ws.Range("J582").Value = 60
$ws.Range("K582").Value = 28000
$ws.Range("L582").Value = 28000
$ws.Range("M582").Value = 28000
$ws.Range("N582").Value = "`$/caja 15 kilos"
$ws.Range("O582").Value = "Región de Arica y Parinacota"
$ws.Range("P582").Value = 1867
$ws.Range("Q582").Value = 15
$ws.Range("R582").Value = "Hortaliza"

# New row 583: Zafiro verde, Primera, serial 45265
$ws.Range("A583").Value = 7
$ws.Range("B583").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C583").Value = "Ñuble"
$ws.Range("D583").Value = 45265
$ws.Range("E583").Value = 16
$ws.Range("F583").Value = 100112002
$ws.Range("G583").Value = "Pimiento"
$ws.Range("H583").Value = "Zafiro verde"
$ws.Range("I583").Value = "Primera"
$ws.Range("J583").Value = 60
$ws.Range("K583").Value = 23000
$ws.Range("L583").Value = 23000
$ws.Range("M583").Value = 23000
$ws.Range("N583").Value = "`$/caja 15 kilos"
$ws.Range("O583").Value = "Región de Arica y Parinacota"
$ws.Range("P583").Value = 1533
$ws.Range("Q583").Value = 15
$ws.Range("R583").Value = "Hortaliza"
